# Update "想去人数" (want-to-go count) values in column F for rows 2, 6, 8, 9, 10
# on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    2  = 588
    6  = 39
    8  = 538
    9  = 3686
    10 = 64
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
